# [Excel] (Shape) Add get active shape image snippet
#
# The "Snippets" table gains one new row (Workbook.getActiveShapeOrNullObject),
# inserted alphabetically between the existing "Workbook.getActiveCell" and
# "Workbook.getSelectedRanges" rows (worksheet row 336, pushing the former
# row 336 and everything below it down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$insertAt = 336

# Push every existing row at/after the insertion point down by one, then
# populate the freed-up row with the new snippet metadata.
$ws.Rows.Item($insertAt).Insert()

$ws.Range("A" + $insertAt).Value = "Excel"
$ws.Range("B" + $insertAt).Value = "Workbook"
$ws.Range("C" + $insertAt).Value = "getActiveShapeOrNullObject"
$ws.Range("D" + $insertAt).Value = 1
$ws.Range("E" + $insertAt).Value = "excel-shape-get-active"
$ws.Range("F" + $insertAt).Value = "getActiveShape"

# Grow the table/autofilter definition so the new last row (390) is included.
$lastRow = $ws.UsedRange.Rows.Count
$lo.Resize($ws.Range("A1:F" + $lastRow))

# Restore the active selection to mirror where editing left off.
$ws.Range("F337").Select()
